$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.924.06'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').Value = '1.552.74'
$ws.Range('E3').Value = '  +1.40%  '
$ws.Range('E4').Value = '  +0.52%  '
$ws.Range('D5').Value = "'206.78"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('E6').Value = '  +0.61%  '
$ws.Range('E7').Value = '  +0.51%  '
$ws.Range('E8').Value = '  +1.70%  '
$ws.Range('E9').Value = '  +1.48%  '
$ws.Range('E10').Value = '  +1.15%  '
$ws.Range('D11').Value = "'0.0858"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.71%  '
$ws.Range('D12').Value = '1.773.06'
$ws.Range('E12').Value = '  +1.37%  '
$ws.Range('D13').Value = '1.551.40'
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('E14').Value = '  +1.59%  '
$ws.Range('D15').Value = "'0.516"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.01%  '
$ws.Range('D16').Value = '26.912.05'
$ws.Range('E16').Value = '  +0.84%  '
$ws.Range('D17').Value = "'61.66"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.99%  '
$ws.Range('D18').Value = "'216.87"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.21%  '
$ws.Range('E19').Value = '  +1.04%  '
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('E22').Value = '  +1.08%  '
$ws.Range('D23').Value = "'9.23"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E24').Value = '  +0.55%  '
$ws.Range('D25').Value = "'153.78"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.27%  '
$ws.Range('D26').Value = "'6.59"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.16%  '
$ws.Range('D27').Value = "'14.87"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.53%  '
$ws.Range('E28').Value = '  +0.45%  '
$ws.Range('E29').Value = '  +1.25%  '
$ws.Range('E30').Value = '  +3.09%  '
$ws.Range('E31').Value = '  +0.13%  '
$ws.Range('D32').Value = "'3.22"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.55%  '
$ws.Range('D33').Value = '1.420.79'
$ws.Range('E33').Value = '  +4.46%  '
$ws.Range('E34').Value = '  +3.08%  '
$ws.Range('E35').Value = '  +4.11%  '
$ws.Range('E36').Value = '  +1.75%  '
$ws.Range('E37').Value = '  +1.01%  '
$ws.Range('E38').Value = '  +0.78%  '
$ws.Range('E39').Value = '  +0.09%  '
$ws.Range('E40').Value = '  +1.58%  '
$ws.Range('E41').Value = '  +0.52%  '
$ws.Range('D42').Value = "'5.70"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('E43').Value = '  -0.48%  '
$ws.Range('E44').Value = '  +2.89%  '
$ws.Range('D45').Value = "'63.70"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.04%  '
$ws.Range('E46').Value = '  -0.13%  '
$ws.Range('D47').Value = '1.687.97'
$ws.Range('E47').Value = '  +1.43%  '
$ws.Range('D48').Value = "'86.23"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.10%  '
$ws.Range('E49').Value = '  +4.24%  '
$ws.Range('D50').Value = '0.0₆0101'
$ws.Range('E50').Value = '  +4.01%  '
$ws.Range('E51').Value = '  +1.68%  '
